$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.737.81'
$ws.Range("E2").Value = '  -2.57%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.394.05'
$ws.Range("E3").Value = '  -3.50%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.93'
$ws.Range("E5").Value = '  -3.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.47'
$ws.Range("E6").Value = '  -6.90%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.395.77'
$ws.Range("E8").Value = '  -3.46%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.473'
$ws.Range("E9").Value = '  -3.07%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.35'
$ws.Range("E10").Value = '  -3.82%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.119'
$ws.Range("E11").Value = '  -4.19%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.374'
$ws.Range("E12").Value = '  -3.72%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.984.77'
$ws.Range("E13").Value = '  -3.19%  '

# Row 14
$ws.Range("E14").Value = '  -0.69%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.404.68'
$ws.Range("E15").Value = '  -3.18%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").Value = '  -5.68%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.778.21'
$ws.Range("E17").Value = '  -2.49%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.56'
$ws.Range("E18").Value = '  -5.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.45'
$ws.Range("E19").Value = '  -5.53%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.62'
$ws.Range("E20").Value = '  -2.65%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.04'
$ws.Range("E21").Value = '  -4.13%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '372.93'
$ws.Range("E22").Value = '  -5.50%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.554'
$ws.Range("E23").Value = '  -3.93%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.536.72'
$ws.Range("E24").Value = '  -3.29%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.21%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.70'
$ws.Range("E26").Value = '  -4.04%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000106'
$ws.Range("E27").Value = '  -10.04%  '

# Row 28
$ws.Range("E28").Value = '  +0.19%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.92'
$ws.Range("E29").Value = '  -6.96%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.14'
$ws.Range("E30").Value = '  -5.72%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.73'
$ws.Range("E31").Value = '  -7.09%  '

# Row 32
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.02%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.149'
$ws.Range("E33").Value = '  -4.86%  '

# Row 34
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.38'
$ws.Range("E34").Value = '  -6.18%  '

# Row 35
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.429.56'
$ws.Range("E35").Value = '  -3.33%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.68'
$ws.Range("E36").Value = '  -3.27%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.19'
$ws.Range("E37").Value = '  -2.88%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '166.00'
$ws.Range("E38").Value = '  -0.68%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.64'
$ws.Range("E39").Value = '  -4.95%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.47'
$ws.Range("E40").Value = '  -5.47%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0750'
$ws.Range("E41").Value = '  -5.11%  '

# Row 42
$ws.Range("E42").Value = '  +0.08%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.773'
$ws.Range("E43").Value = '  -4.87%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.48'
$ws.Range("E44").Value = '  -1.90%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.24'
$ws.Range("E45").Value = '  -4.78%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.74'
$ws.Range("E46").Value = '  -10.42%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.54'
$ws.Range("E47").Value = '  -7.16%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").Value = '  -9.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.62'
$ws.Range("E49").Value = '  -3.20%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.258.35'
$ws.Range("E50").Value = '  -5.56%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.849'
$ws.Range("E51").Value = '  -5.61%  '
